$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("About")
$ws2 = $wb.Worksheets.Item("FPCbS")

# ------------------------------------------------------------------
# FPCbS sheet: header row rework
# ------------------------------------------------------------------
# B1 used to just say "FP"; it now carries a fuller label, is bold,
# wraps, and sits next to a (blank) styled A1 header cell.
$ws2.Range("B1").WrapText = $true
$ws2.Range("A1").Font.Bold = $true
$ws2.Range("A1").WrapText = $true

# New rows of fuel sources with flexibility points of 0
$ws2.Range("A15").Value = "crude oil"
$ws2.Range("B15").Value = 0
$ws2.Range("A16").Value = "heavy or residual fuel oil"
$ws2.Range("B16").Value = 0
$ws2.Range("A17").Value = "municipal solid waste"
$ws2.Range("B17").Value = 0

$ws2.Range("B1").Value = "FP (flexibility points/MW)"

$ws2.Rows.Item(1).RowHeight = 75

$ws2.Columns.Item(1).ColumnWidth = 25
$ws2.Columns.Item(2).ColumnWidth = 13.665

$ws2.Activate()
$ws2.Range("B1").Select()

# ------------------------------------------------------------------
# About sheet: add a new note for India regional adaptation
# ------------------------------------------------------------------
$ws1.Range("A46").Value = "For India - solar PV, onshore and offshore wind need flexibility."

$ws1.Activate()
$ws1.Range("A46").Select()
